# Adds a new "2022-Q1" sheet (positioned between "2021-Q4" and "总计") that
# holds the full 2022-Q1 fund-holdings table, and updates the "总计" (totals)
# summary sheet with a new leading row for 2022-Q1 (pushing the existing
# 2021-Q4 summary row down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" worksheet right after "2021-Q4".
# ---------------------------------------------------------------------
$src = $wb.Worksheets.Item("2021-Q4")
$new = $wb.Worksheets.Add($null, $src)
$new.Name = "2022-Q1"

# Clone the header-row formatting (bold + border, style used on "2021-Q4"!B1:H1)
# and the formatting used on the "2021-Q4"!A column (style used for the
# row-sequence column) so the new sheet visually matches its sibling.
$src.Range("B1:H1").Copy()
$new.Range("B1:H1").PasteSpecial(-4122)

$src.Range("A2").Copy()
$new.Range("A2:A18").PasteSpecial(-4122)

# Columns B,D,E,F,G hold numeric-looking text values (fund code / size /
# position percentages kept as strings, exactly like on "2021-Q4") -- force
# text format BEFORE writing so those strings (e.g. "000971", "15.43") are
# not coerced into numbers. Column C (fund name) is never numeric-looking so
# it needs no special handling.
$new.Range("B2:B18").NumberFormat = "@"
$new.Range("D2:G18").NumberFormat = "@"

$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

$rows = @(
    @(0, "000971", "诺安新经济股票", "15.43", "82.95", "8.11", "1.2514", 2),
    @(1, "002446", "广发利鑫灵活配置混合A", "12.46", "74.35", "5.15", "0.6417", 4),
    @(2, "008328", "诺安新兴产业混合", "4.83", "68.49", "8.14", "0.3932", 3),
    @(3, "200012", "长城中小盘成长混合", "12.65", "84.26", "2.25", "0.2846", 7),
    @(4, "001127", "中银宏观策略灵活配置混合", "4.38", "74.26", "3.41", "0.1494", 4),
    @(5, "010602", "长城均衡优选混合", "4.55", "84.61", "2.36", "0.1074", 6),
    @(6, "002291", "诺安安鑫灵活配置混合", "2.19", "81.55", "3.53", "0.0773", 8),
    @(7, "011172", "广发利鑫灵活配置混合C", "1.10", "74.35", "5.15", "0.0566", 4),
    @(8, "000826", "广发中证百度百发策略100指数A", "4.11", "92.42", "1.15", "0.0473", 2),
    @(9, "000827", "广发中证百度百发策略100指数E", "4.11", "92.42", "1.15", "0.0473", 2),
    @(10, "001780", "诺安改革趋势灵活配置混合", "0.46", "68.34", "9.76", "0.0449", 1),
    @(11, "002137", "诺安利鑫灵活配置混合", "0.45", "76.34", "4.43", "0.0199", 2),
    @(12, "000591", "中银健康生活混合", "0.55", "70.29", "3.39", "0.0186", 3),
    @(13, "000120", "中银美丽中国混合", "0.50", "81.97", "3.52", "0.0176", 4),
    @(14, "005545", "中银改革红利灵活配置混合", "0.51", "68.74", "3.01", "0.0154", 5),
    @(15, "003717", "中银量化精选灵活配置混合A", "0.49", "90.38", "1.14", "0.0056", 10),
    @(16, "010484", "中银量化精选灵活配置混合C", "0.01", "90.38", "1.14", "0.0001", 10)
)

$r = 2
foreach ($row in $rows) {
    $new.Cells.Item($r, 1).Value = $row[0]
    $new.Cells.Item($r, 2).Value = $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = $row[3]
    $new.Cells.Item($r, 5).Value = $row[4]
    $new.Cells.Item($r, 6).Value = $row[5]
    $new.Cells.Item($r, 7).Value = $row[6]
    $new.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# The values are now stored as genuine text (Excel picked the type up from
# the "@" format at entry time); drop that temporary formatting again so the
# cells end up plain/un-styled, exactly like the matching text cells on
# "2021-Q4".
$new.Range("B2:B18").ClearFormats()
$new.Range("D2:G18").ClearFormats()

# ---------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new summary row for
#    2022-Q1 above the existing 2021-Q4 row.
# ---------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Give the new A3 row-index cell the same style as A2 (bordered / centered)
# before moving the 2021-Q4 figures down into it.
$tot.Range("A2").Copy()
$tot.Range("A3").PasteSpecial(-4122)

$tot.Range("A3").Value = 1
$tot.Range("B3").Value = "2021-Q4"
$tot.Range("C3").Value = 5
$tot.Range("D3").Value = 0.24

$tot.Range("A2").Value = 0
$tot.Range("B2").Value = "2022-Q1"
$tot.Range("C2").Value = 17
$tot.Range("D2").Value = 3.18

Write-Output "2022-Q1 sheet added; 总计 sheet updated"
